$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the value of B9 to the new text "みなさん" (keeps existing style s="3")
$ws.Range("B9").Value = "みなさん"

# Update the selection to B9 (matches sheetView selection change in the diff)
$ws.Range("B9").Select()
